$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 115435.63
$ws.Range("I132").Value = 1586.9851
$ws.Range("J132").Value = 478667.06
$ws.Range("K132").Value = 4760.955300000001
$ws.Range("L132").Value = 1436001.18
$ws.Range("M132").Value = -2230.955300000001
$ws.Range("N132").Value = -1441061.18
# Row 135
$ws.Range("H135").Value = 9616340
$ws.Range("I135").Value = 275.57574
$ws.Range("J135").Value = 26317926
$ws.Range("K135").Value = 2480.18166
$ws.Range("L135").Value = 236861334
$ws.Range("M135").Value = 54.81833999999981
$ws.Range("N135").Value = -236866404
# Row 137
$ws.Range("H137").Value = 28513.074
$ws.Range("I137").Value = 44427.176
$ws.Range("J137").Value = 6982.2354
$ws.Range("K137").Value = 133281.528
$ws.Range("L137").Value = 20946.7062
$ws.Range("M137").Value = -130731.528
$ws.Range("N137").Value = -26046.7062
# Row 138
$ws.Range("H138").Value = 1606.56
$ws.Range("I138").Value = 824.2
$ws.Range("J138").Value = 2246.6726
$ws.Range("K138").Value = 2472.6
$ws.Range("L138").Value = 6740.0178
$ws.Range("M138").Value = 2667.4
$ws.Range("N138").Value = -17020.0178

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2337.44
$ws.Range("I32").Value = 1959.427
$ws.Range("J32").Value = 5395.909
$ws.Range("K32").Value = 1959.427
$ws.Range("L32").Value = 5395.909
$ws.Range("M32").Value = -1672.427
$ws.Range("N32").Value = -5969.909
# Row 45
$ws.Range("H45").Value = 956.129
$ws.Range("I45").Value = 976.1053000000001
$ws.Range("J45").Value = 924.5
$ws.Range("K45").Value = 976.1053000000001
$ws.Range("L45").Value = 924.5
$ws.Range("M45").Value = -599.1053000000001
$ws.Range("N45").Value = -1678.5
# Row 61
$ws.Range("H61").Value = 1013.8333
$ws.Range("I61").Value = 988.9259
$ws.Range("J61").Value = 1238
$ws.Range("K61").Value = 988.9259
$ws.Range("L61").Value = 1238
$ws.Range("M61").Value = -776.9259
$ws.Range("N61").Value = -1662
# Row 74
$ws.Range("H74").Value = 16936.523
$ws.Range("I74").Value = 22709.305
$ws.Range("J74").Value = 1316.0588
$ws.Range("K74").Value = 22709.305
$ws.Range("L74").Value = 1316.0588
$ws.Range("M74").Value = -21835.305
$ws.Range("N74").Value = -3064.0588
# Row 77
$ws.Range("H77").Value = 16936.523
$ws.Range("I77").Value = 22709.305
$ws.Range("J77").Value = 1316.0588
$ws.Range("K77").Value = 113546.525
$ws.Range("L77").Value = 6580.294
$ws.Range("M77").Value = -109178.525
$ws.Range("N77").Value = -15316.294
# Row 136
$ws.Range("H136").Value = 1013.8333
$ws.Range("I136").Value = 988.9259
$ws.Range("J136").Value = 1238
$ws.Range("K136").Value = 2966.7777
$ws.Range("L136").Value = 3714
$ws.Range("M136").Value = -416.7776999999996
$ws.Range("N136").Value = -8814

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 15116.487
$ws.Range("I134").Value = 834.7846
$ws.Range("J134").Value = 86525
$ws.Range("K134").Value = 2504.3538
$ws.Range("L134").Value = 259575
$ws.Range("M134").Value = 30.64620000000014
$ws.Range("N134").Value = -264645

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 23847.193
$ws.Range("I31").Value = 25838.56
$ws.Range("J31").Value = 15549.833
$ws.Range("K31").Value = 25838.56
$ws.Range("L31").Value = 15549.833
$ws.Range("M31").Value = -25543.56
$ws.Range("N31").Value = -16139.833
# Row 34
$ws.Range("H34").Value = 23847.193
$ws.Range("I34").Value = 25838.56
$ws.Range("J34").Value = 15549.833
$ws.Range("K34").Value = 25838.56
$ws.Range("L34").Value = 15549.833
$ws.Range("M34").Value = -25636.56
$ws.Range("N34").Value = -15953.833
# Row 58
$ws.Range("H58").Value = 849.9394
$ws.Range("I58").Value = 545.2195
$ws.Range("J58").Value = 1349.68
$ws.Range("K58").Value = 545.2195
$ws.Range("L58").Value = 1349.68
$ws.Range("M58").Value = -342.2195
$ws.Range("N58").Value = -1755.68
# Row 134
$ws.Range("H134").Value = 831.68494
$ws.Range("I134").Value = 772.7347
$ws.Range("J134").Value = 952.0417
$ws.Range("K134").Value = 2318.2041
$ws.Range("L134").Value = 2856.1251
$ws.Range("M134").Value = 216.7959000000001
$ws.Range("N134").Value = -7926.1251
# Row 136
$ws.Range("H136").Value = 849.9394
$ws.Range("I136").Value = 545.2195
$ws.Range("J136").Value = 1349.68
$ws.Range("K136").Value = 1635.6585
$ws.Range("L136").Value = 4049.04
$ws.Range("M136").Value = 914.3415
$ws.Range("N136").Value = -9149.040000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 3452.9
$ws.Range("I75").Value = 3333
$ws.Range("J75").Value = 3459.2104
$ws.Range("K75").Value = 9999
$ws.Range("L75").Value = 10377.6312
$ws.Range("M75").Value = -9001
$ws.Range("N75").Value = -12373.6312
# Row 78
$ws.Range("H78").Value = 3452.9
$ws.Range("I78").Value = 3333
$ws.Range("J78").Value = 3459.2104
$ws.Range("K78").Value = 29997
$ws.Range("L78").Value = 31132.8936
$ws.Range("M78").Value = -25005
$ws.Range("N78").Value = -41116.8936

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
# Row 122
$ws.Range("H122").Value = 1090.2858
$ws.Range("I122").Value = 977.8182
$ws.Range("K122").Value = 2933.4546
$ws.Range("M122").Value = -483.4546

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2659.9167
$ws.Range("I61").Value = 2371.2856
$ws.Range("J61").Value = 3064
$ws.Range("K61").Value = 2371.2856
$ws.Range("L61").Value = 3064
$ws.Range("M61").Value = -2169.2856
$ws.Range("N61").Value = -3468
# Row 113
$ws.Range("H113").Value = 2659.9167
$ws.Range("I113").Value = 2371.2856
$ws.Range("J113").Value = 3064
$ws.Range("K113").Value = 2371.2856
$ws.Range("L113").Value = 3064
$ws.Range("M113").Value = -201.2856000000002
$ws.Range("N113").Value = -7404
# Row 132
$ws.Range("H132").Value = 196166.42
$ws.Range("I132").Value = 41605.88
$ws.Range("K132").Value = 124817.64
$ws.Range("M132").Value = -122287.64
# Row 136
$ws.Range("H136").Value = 201392.52
$ws.Range("I136").Value = 295199.9
$ws.Range("J136").Value = 2051.8125
$ws.Range("K136").Value = 885599.7000000001
$ws.Range("L136").Value = 6155.4375
$ws.Range("M136").Value = -883049.7000000001
$ws.Range("N136").Value = -11255.4375

$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Range("H114").Value = 30420
$ws.Range("J114").Value = 30420
$ws.Range("L114").Value = 30420
$ws.Range("N114").Value = -39098
# Row 122
$ws.Range("H122").Value = 7778.5713
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
# Row 132
$ws.Range("H132").Value = 2505.7896
$ws.Range("I132").Value = 666.0263
$ws.Range("J132").Value = 6185.316
$ws.Range("K132").Value = 1998.0789
$ws.Range("L132").Value = 18555.948
$ws.Range("M132").Value = 531.9211
$ws.Range("N132").Value = -23615.948
# Row 136
$ws.Range("H136").Value = 1198500
$ws.Range("I136").Value = 1458796.4
$ws.Range("J136").Value = 527209.4399999999
$ws.Range("K136").Value = 4376389.199999999
$ws.Range("L136").Value = 1581628.32
$ws.Range("M136").Value = -4373839.199999999
$ws.Range("N136").Value = -1586728.32
